$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for caseid=1087 / category="ECG Lead II" / tname="SNUADC/ECG_II" /
# tid="0ed1f6a51b46c3f6436e61b2abc69079eb9bcf2c" (sheet row 34) was removed.
# Deleting it shifts every following row up by one, which is exactly what the
# diff shows (dimension A1:D237 -> A1:D236, all subsequent rows renumbered).
$ws.Rows(34).Delete()

# The data rows (everything below the header) got a bit taller.
$lastRow = $ws.Cells(1, 1).End(-4121).Row
$dataRange = $ws.Range(("A2:D{0}" -f $lastRow))
$dataRange.RowHeight = 19.5

# The caseid column's font was switched from the theme text color to an
# explicit black RGB color.
$caseIdRange = $ws.Range(("A2:A{0}" -f $lastRow))
$caseIdRange.Font.Color = 0
